$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title heading and bold repeated text (Find/Replace with wdReplaceAll covers both occurrences)
Replace-Text "Play Esqueleto Explosivo - Free Slot Game Review" "Play Esqueleto Explosivo for Free"

# "What we like" bullet list
Replace-Text "Captivating and colorful 3D graphics" "Explosive and entertaining gameplay"
Replace-Text "Smooth and entertaining gameplay" "High-quality 3D graphics and animations"
Replace-Text "Exclusive special characters increase winning potential" "Unique and rewarding special characters"
Replace-Text "Mexican-style music and theme" "Captivating Mexican-themed setting"

# "What we don't like" bullet list
Replace-Text "No progressive jackpot" "Lack of bonus features"

# Meta description (italic)
Replace-Text "Read our review of Esqueleto Explosivo, a free online slot game. Explore its captivating 3D graphics and exclusive special characters that increase winning potential." "Discover the explosive gameplay and captivating graphics of Esqueleto Explosivo in this free slot review."
